$d = $word.ActiveDocument

# The document-number line currently reads "Số: (1)/24/TTr-KXN". The
# requested edit bumps the two-digit year from "24" to "25" -- i.e. a
# user selected just the single digit "4" and retyped it as "5". Word
# keeps that kind of in-place single-character edit as its own run
# rather than re-merging it with the untouched text before/after it,
# so the run ends up split in three: ".../2", "5", "/TTr-KXN...".

$r = $d.Content
$found = $r.Find.Execute("Số: (1)/24/TTr-KXN")

if ($found) {
    $matchStart = $r.Start
    $matchEnd = $r.End

    # Bookmark the whole matched run first. This pins its left edge so
    # that when we start rewriting text inside it, the engine doesn't
    # fold the edited run back together with the separate (unrelated)
    # run holding the leading whitespace that precedes it.
    $wholeRange = $d.Range($matchStart, $matchEnd)
    $guardBookmark = $d.Bookmarks.Add("zzEditGuard", $wholeRange)

    # Find "24" within the matched text, then narrow down to just its
    # second character (the "4") -- that's the only character that
    # actually changes.
    $numRange = $d.Range($matchStart, $matchEnd)
    $numFound = $numRange.Find.Execute("24")
    $digitRange = $d.Range($numRange.Start + 1, $numRange.End)
    $digitRange.Text = "5"

    # Bookmark just the replaced "5" momentarily. This keeps it as an
    # isolated run instead of letting it re-coalesce with its
    # neighbours once we touch it, matching how Word leaves a manual
    # single-character retype as its own run. The bookmarks themselves
    # are not part of the intended change, so drop them right away.
    $splitRange = $d.Range($numRange.Start + 1, $numRange.Start + 2)
    $splitBookmark = $d.Bookmarks.Add("zzEditSplit", $splitRange)

    $d.Bookmarks("zzEditSplit").Delete()
    $d.Bookmarks("zzEditGuard").Delete()
}
